$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in column A (rows 10, 12, 18) and column C (row 13)
$ws.Range("A10").Value = -20.926
$ws.Range("A12").Value = -21.694
$ws.Range("C13").Value = -13.059
$ws.Range("A18").Value = -21.694
